# "Generate Report for Handback"
#
# For each locale sheet (zh-cn, de-de):
#  - Status column (C) moves from "Ready for handoff" to
#    "Handed back: in sync with en-US" (also mirrored on the Overview
#    sheet, since it shares the same status text).
#  - The "Latest Target File" (F) / "Latest Handback File" (G) columns
#    are populated with the same handoff .md / .xlf links that are
#    already present in columns A / D, now that the package has been
#    handed back.
#  - The "Latest Handback DateTime" (H) is stamped with the real
#    handback time instead of the zero-date placeholder.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Update every cell that currently shows the old status text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: populate F/G (Latest Target File / Latest Handback
#    File) and stamp the Latest Handback DateTime (H).
# ---------------------------------------------------------------------
$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/6275a4da9509d6311e3af691546b4979c75e8a65/e2e/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.md",
    [Type]::Missing,
    [Type]::Missing,
    "2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.md")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2c4d779ba04829f60b44c18ecce367a5e341734c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.ccee5f170c90ec7c63c04517415a7f4e04a48849.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.ccee5f170c90ec7c63c04517415a7f4e04a48849.zh-cn.xlf")

$wsZhCn.Range("H2").Value = "2016-03-21 15:02:36"

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/6275a4da9509d6311e3af691546b4979c75e8a65/e2e/e2d5556a-efb3-4967-a0ea-247ca6604ea2.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2d5556a-efb3-4967-a0ea-247ca6604ea2.md")

$wsZhCn.Hyperlinks.Add(
    $wsZhCn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2c4d779ba04829f60b44c18ecce367a5e341734c/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/e2d5556a-efb3-4967-a0ea-247ca6604ea2.4fc15d15cc924eec3f2e0ddbd167bb8ca4a7aeb0.zh-cn.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "e2d5556a-efb3-4967-a0ea-247ca6604ea2.4fc15d15cc924eec3f2e0ddbd167bb8ca4a7aeb0.zh-cn.xlf")

$wsZhCn.Range("H3").Value = "2016-03-21 15:02:36"

# ---------------------------------------------------------------------
# 3. de-de sheet: same treatment, with its own handback timestamp.
# ---------------------------------------------------------------------
$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/6275a4da9509d6311e3af691546b4979c75e8a65/e2e/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.md",
    [Type]::Missing,
    [Type]::Missing,
    "2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.md")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f87811625197db4aa2d219e29072069f73769b2e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.ccee5f170c90ec7c63c04517415a7f4e04a48849.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "2d7c8c26-fe6a-4d8b-88c7-3faa8d4c292c.ccee5f170c90ec7c63c04517415a7f4e04a48849.de-de.xlf")

$wsDeDe.Range("H2").Value = "2016-03-21 15:02:43"

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/6275a4da9509d6311e3af691546b4979c75e8a65/e2e/e2d5556a-efb3-4967-a0ea-247ca6604ea2.md",
    [Type]::Missing,
    [Type]::Missing,
    "e2d5556a-efb3-4967-a0ea-247ca6604ea2.md")

$wsDeDe.Hyperlinks.Add(
    $wsDeDe.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/f87811625197db4aa2d219e29072069f73769b2e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/e2d5556a-efb3-4967-a0ea-247ca6604ea2.4fc15d15cc924eec3f2e0ddbd167bb8ca4a7aeb0.de-de.xlf",
    [Type]::Missing,
    [Type]::Missing,
    "e2d5556a-efb3-4967-a0ea-247ca6604ea2.4fc15d15cc924eec3f2e0ddbd167bb8ca4a7aeb0.de-de.xlf")

$wsDeDe.Range("H3").Value = "2016-03-21 15:02:43"

Write-Host "Handback report generated."
